$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the incorrectly-included "The Merchant of Venice" row (row 23).
#    Excel shifts every following row up by one and the now-unused
#    shared string is dropped automatically on save.
$ws.Rows("23").Delete()

# 2. The previous row's Vietnamese-translation flag ("Y") was wrong for
#    every other William Shakespeare title - correct each one to "N".
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $author = $ws.Cells($r, 2).Value2
    if ($author -eq "William Shakespeare") {
        $ws.Cells($r, 3).Value = "N"
    }
}

# 3. Re-point the autofilter / filter-database range to the new extent
#    of the table (it used to stop at row 22, well short of the data).
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:D81").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$81"
    }
}

# 4. Leave the selection on the row that now occupies position 23
#    (mirrors the saved view state) instead of the old scrolled-down spot.
$ws.Rows("23").Select()
